$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1088
$ws.Range("F7").Value = 134
$ws.Range("F8").Value = 1390
$ws.Range("F14").Value = 1295
$ws.Range("F15").Value = 435
$ws.Range("F16").Value = 460
$ws.Range("F19").Value = 604
$ws.Range("F20").Value = 2527
$ws.Range("F29").Value = 352
$ws.Range("F31").Value = 44
$ws.Range("F32").Value = 28

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 590
$ws.Range("F6").Value = 590

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1750
$ws.Range("F5").Value = 2042
$ws.Range("F10").Value = 1082

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1750
$ws.Range("F5").Value = 2042
$ws.Range("F12").Value = 1082
$ws.Range("F17").Value = 1088
$ws.Range("F18").Value = 134
$ws.Range("F19").Value = 590
$ws.Range("F25").Value = 1295
$ws.Range("F26").Value = 435
$ws.Range("F27").Value = 460
$ws.Range("F30").Value = 605
$ws.Range("F38").Value = 352
$ws.Range("F41").Value = 44
